$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 140 -- shifts existing rows 140:223 down to 141:224
$ws.Rows.Item(140).EntireRow.Insert()

# Populate the newly inserted row 140 with the new record's data.
# (A, B, C, E, F, G, H, I, N, O, Q, R mirror the rest of this table's rows.)
$ws.Range("A140").Value = 4
$ws.Range("B140").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C140").Value = "Los Lagos"
$ws.Range("D140").Value = 44582
$ws.Range("E140").Value = 10
$ws.Range("F140").Value = 100112037
$ws.Range("G140").Value = "Cebollín"
$ws.Range("H140").Value = "Sin especificar"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 180
$ws.Range("K140").Value = 6000
$ws.Range("L140").Value = 6000
$ws.Range("M140").Value = 6000
$ws.Range("N140").Value = "`$/paquete 36 unidades"
$ws.Range("O140").Value = "Región Metropolitana"
$ws.Range("P140").Value = 167
$ws.Range("Q140").Value = 36
$ws.Range("R140").Value = "Hortaliza"

# Match the date-formatted style already used by the other rows' "Fecha" column.
$ws.Range("D140").NumberFormat = $ws.Range("D141").NumberFormat
